$d = $word.ActiveDocument
$d.Content.Find.Execute("Suhendar ()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Suhendar (1127050153)", 2)
